# Applies the "Add data for 2022-01-07" update:
#  - Rename sheet / update shared-string header from "December 29" to "December 30"
#  - Bump a handful of December monthly carjacking counts (current + prior years)
#    to reflect the newly-added day of data (Dec 30, 2021).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet name ---
$ws.Name = "Through 2021-12-30"

# --- Header cell text (shared string used by B1) ---
$ws.Range("B1").Value = "December 2021 (through December 30)"

# --- Cell value updates / additions ---
# Row 2
$ws.Range("BV2").Value = 2

# Row 4
$ws.Range("N4").Value = 14

# Row 6
$ws.Range("N6").Value = 12

# Row 7
$ws.Range("AX7").Value = 10
$ws.Range("BJ7").Value = 7
$ws.Range("BV7").Value = 7

# Row 8
$ws.Range("B8").Value = 7
$ws.Range("AL8").Value = 2
$ws.Range("BJ8").Value = 7

# Row 9
$ws.Range("B9").Value = 8

# Row 10
$ws.Range("Z10").Value = 2

# Row 18
$ws.Range("Z18").Value = 2

# Row 21
$ws.Range("B21").Value = 5
$ws.Range("AL21").Value = 2

# Row 22 (new cell)
$ws.Range("B22").Value = 1

# Row 23
$ws.Range("AL23").Value = 5

# Row 27
$ws.Range("B27").Value = 3
$ws.Range("N27").Value = 1

# Row 28 (new cell)
$ws.Range("B28").Value = 1

# Row 30
$ws.Range("B30").Value = 5

# Row 33
$ws.Range("B33").Value = 5

# Row 34
$ws.Range("N34").Value = 3

# Row 40
$ws.Range("N40").Value = 6

# Row 41
$ws.Range("N41").Value = 3

# Row 49
$ws.Range("B49").Value = 5
$ws.Range("AX49").Value = 1

# Row 53
$ws.Range("N53").Value = 2

# Row 54
$ws.Range("N54").Value = 3

# Row 55 (new cell)
$ws.Range("N55").Value = 1

# Row 61 (new cell)
$ws.Range("AX61").Value = 1

# Row 66
$ws.Range("BJ66").Value = 2

# Row 79
$ws.Range("B79").Value = 2
